# Remove the "Dipus" and "Eothenomys" genus rows from the blood-meal
# detection table (one pair of rows per sample block: K1, K2, K3, K4, K5,
# M1, V2, V3). Deleting shifts subsequent rows up, so delete from the
# bottom of the sheet upward to keep row numbers stable while iterating.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(10,11,20,21,30,31,40,41,50,51,60,61,70,71,80,81)
$rowsToDelete = $rowsToDelete | Sort-Object -Descending

foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

$ws.Range("F65").Select()
